$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing revision (column B) values to reflect re-based GVA figures
$ws.Range("B7").Value = 19.58999999999999
$ws.Range("B23").Value = 20.8
$ws.Range("B30").Value = 18.46000000000001
$ws.Range("B40").Value = 10.31
$ws.Range("B44").Value = 9.599999999999994
$ws.Range("B45").Value = 10.21000000000001
$ws.Range("B48").Value = 9.560000000000002
$ws.Range("B52").Value = 10.89
$ws.Range("B53").Value = 11.27
$ws.Range("B55").Value = 10.99000000000001
$ws.Range("B56").Value = 11.17999999999999
$ws.Range("B59").Value = 2.359999999999999
$ws.Range("B68").Value = 0.2199999999999989
$ws.Range("B69").Value = 1.269999999999996
$ws.Range("B72").Value = 1.88000000000001
$ws.Range("B73").Value = 2.02000000000001
$ws.Range("B75").Value = 1.109999999999999
$ws.Range("B76").Value = 2.359999999999999
$ws.Range("B77").Value = 0.7900000000000063
$ws.Range("B80").Value = 0.4299999999999926
$ws.Range("B81").Value = -0.06000000000000227

# Append the new quarterly datapoint as row 82, matching the date-cell style of A81
$ws.Cells.Item(81, 1).Copy()
$ws.Cells.Item(82, 1).PasteSpecial(-4122)
$ws.Cells.Item(82, 1).Value = 45884
$ws.Cells.Item(82, 2).Value = 0
